$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (not auto-numeric) formatting for Price cells whose new value
# would otherwise be auto-parsed as a number by Excel.
$textCells = @("D5","D6","D8","D9","D10","D11","D14","D15","D18","D20","D21","D22","D23","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D41","D42","D44","D45","D46","D49","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated Price (D) / Volume(1h) (E) values row by row.
$ws.Range("D2").Value = "29.449.61"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "1.850.84"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "240.96"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").Value = "0.6303"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.07703"
$ws.Range("E8").Value = "  +2.49%  "
$ws.Range("D9").Value = "0.2930"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").Value = "24.76"
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("D11").Value = "0.07743"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("D12").Value = "1.880.09"
$ws.Range("E12").Value = "  +1.85%  "
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").Value = "0.00001075"
$ws.Range("E14").Value = "  +5.14%  "
$ws.Range("D15").Value = "0.6795"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("E16").Value = "  +0.97%  "
$ws.Range("D17").Value = "2.150.48"
$ws.Range("E17").Value = "  +2.51%  "
$ws.Range("D18").Value = "6.212"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").Value = "29.485.41"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").Value = "228.74"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").Value = "12.46"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "7.464"
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "157.51"
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("D26").Value = "0.1382"
$ws.Range("D27").Value = "8.416"
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("D28").Value = "17.71"
$ws.Range("D29").Value = "1.343"
$ws.Range("E29").Value = "  +5.94%  "
$ws.Range("D30").Value = "1.470"
$ws.Range("E30").Value = "  +0.72%  "
$ws.Range("D31").Value = "0.05684"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D32").Value = "4.131"
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("D33").Value = "4.040"
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("D34").Value = "1.852"
$ws.Range("E34").Value = "  +1.29%  "
$ws.Range("D35").Value = "1.164"
$ws.Range("E35").Value = "  +0.98%  "
$ws.Range("D36").Value = "0.7058"
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("D37").Value = "2.585"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").Value = "2.782"
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("E39").Value = "  -0.68%  "
$ws.Range("D40").Value = "1.219.30"
$ws.Range("E40").Value = "  -1.76%  "
$ws.Range("D41").Value = "6.552"
$ws.Range("E41").Value = "  +5.22%  "
$ws.Range("D42").Value = "0.9094"
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").Value = "101.76"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").Value = "66.45"
$ws.Range("E45").Value = "  +1.59%  "
$ws.Range("D46").Value = "0.00000000120"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("E47").Value = "  +0.96%  "
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("D49").Value = "9.007"
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("D50").Value = "1.687"
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("E51").Value = "  +2.41%  "
